$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A and B both to 15.42578125 characters; the
# engine quantizes ColumnWidth to its internal pixel grid, so we pick the
# input that lands on the closest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 14.6

# Update cell values in columns A and B, rows 1-32
$ws.Range("A1").Value = -0.098162608357640124
$ws.Range("B1").Value = 0.097958897956829105
$ws.Range("A2").Value = -0.075854384001919684
$ws.Range("B2").Value = 0.075128502934767027
$ws.Range("A3").Value = -0.025417387845127593
$ws.Range("B3").Value = 0.025299765315587308
$ws.Range("A4").Value = -0.017299765362992048
$ws.Range("B4").Value = 0.016985361334297622
$ws.Range("A5").Value = -0.013985361355779879
$ws.Range("B5").Value = 0.012925872654579784
$ws.Range("A6").Value = 0.0050662406835115803
$ws.Range("B6").Value = -0.0053209592618550516
$ws.Range("A7").Value = 0.015320959201875706
$ws.Range("B7").Value = -0.015378072985277935
$ws.Range("A8").Value = 0.025378072926321771
$ws.Range("B8").Value = -0.025483104696735204
$ws.Range("A9").Value = 0.02748310468277948
$ws.Range("B9").Value = -0.027572336751759519
$ws.Range("A10").Value = 0.029572336740553595
$ws.Range("B10").Value = -0.029577846551406495
$ws.Range("A11").Value = -0.010548763135106221
$ws.Range("B11").Value = 0.010537250822598487
$ws.Range("A12").Value = -0.0070372508417109181
$ws.Range("B12").Value = 0.0069598079325619011
$ws.Range("A13").Value = -0.0034598079525496317
$ws.Range("B13").Value = 0.0034300914012472816
$ws.Range("A14").Value = 0.0045699085547106222
$ws.Range("B14").Value = -0.0045759822632538416
$ws.Range("A15").Value = 0.0055759822566425754
$ws.Range("B15").Value = -0.0055785067032063651
$ws.Range("A16").Value = 0.0075785066916358446
$ws.Range("B16").Value = -0.0075854154381476135
$ws.Range("A17").Value = -0.0040034047542656381
$ws.Range("B17").Value = 0.003999999978631763
$ws.Range("A18").Value = -0.016104513987329483
$ws.Range("B18").Value = 0.016091498877127464
$ws.Range("A19").Value = -0.012091498899027719
$ws.Range("B19").Value = 0.012016661288238417
$ws.Range("A20").Value = -0.00801666131187595
$ws.Range("B20").Value = 0.0080056720098511391
$ws.Range("A21").Value = -0.0040056720337569018
$ws.Range("B21").Value = 0.003999999975943247
$ws.Range("A22").Value = -0.045711115114920986
$ws.Range("B22").Value = 0.045498099295050665
$ws.Range("A23").Value = -0.040498099327407999
$ws.Range("B23").Value = 0.040098833187395932
$ws.Range("A24").Value = -0.020098833302816921
$ws.Range("B24").Value = 0.019999999882991837
$ws.Range("A25").Value = -0.015992113355645543
$ws.Range("B25").Value = 0.015963018488754699
$ws.Range("A26").Value = -0.049197004523461985
$ws.Range("B26").Value = 0.049173917465660466
$ws.Range("A27").Value = -0.046673917486867555
$ws.Range("B27").Value = 0.046538627414024258
$ws.Range("A28").Value = -0.0445386274374302
$ws.Range("B28").Value = 0.044446939443719913
$ws.Range("A29").Value = -0.037446939497132625
$ws.Range("B29").Value = 0.037420617694547076
$ws.Range("A30").Value = -0.021164843965342683
$ws.Range("B30").Value = 0.021022545110775059
$ws.Range("A31").Value = -0.014022545167831524
$ws.Range("B31").Value = 0.01400123360150296
$ws.Range("A32").Value = -0.0040012336748578292
$ws.Range("B32").Value = 0.0039999999584807711
